$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Watermelon"
$ws.Range("B4").Value = 30
$ws.Range("C4").Value = 4

$ws.Range("A5").Value = "Durian"
$ws.Range("B5").Value = 60
$ws.Range("C5").Value = 2

$ws.Range("A2:C3").Copy()
$ws.Range("A4:C5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C5").Select()
